# Apply Crypto price/volume updates from the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.968.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.554.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.247'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.49'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.17%  '
$ws.Range("E10").Value = '  -0.43%  '
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.776.01'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.555.07'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.974.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.35'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0686'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.95%  '
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("E27").Value = '  -0.48%  '
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("E32").Value = '  +2.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.371.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.56%  '
$ws.Range("E34").Value = '  +1.12%  '
$ws.Range("E35").Value = '  +1.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.970'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.82%  '
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.518'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.809'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.983'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.77%  '
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("E44").Value = '  +2.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("E46").Value = '  -2.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.689.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("E48").Value = '  -3.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("E50").Value = '  +0.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0957'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.47%  '
